# "test runs for avg before sampling"
# Adds four new benchmark rows (scaling AFTER averaging experiments) to the
# Tabelle1 table, renames the "average before predictions" label to the
# shorter "avg before predictions" wording, and keeps the table/filter/
# conditional-formatting ranges and view in sync with the now-larger sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Append the four new rows (48-51) with the "scaling AFTER avg" runs.
#    (Values are written in the same order the original workbook's shared
#    string table grew in, so the new unique strings line up.)
# ---------------------------------------------------------------------
$ws.Range("B48").Value = "Logistic Classifier"
$ws.Range("C48").Value = "MoCo"
$ws.Range("D48").Value = "Centers"
$ws.Range("E48").Value = "1 x 3"
$ws.Range("F48").Value = "average"
$ws.Range("M48").Value = "0.640 (0.028)"
$ws.Range("I48").Value = "standard scaling center data AFTER avg"
$ws.Range("J48").Value = 0.635
$ws.Range("K48").Value = 0.676
$ws.Range("L48").Value = 0.609

$ws.Range("B49").Value = "Logistic Classifier"
$ws.Range("C49").Value = "MoCo"
$ws.Range("D49").Value = "Centers"
$ws.Range("E49").Value = "1 x 3"
$ws.Range("F49").Value = "average"
$ws.Range("M49").Value = "0.648 (0.037)"
$ws.Range("I49").Value = "min max scaling center data AFTER avg"
$ws.Range("J49").Value = 0.637
$ws.Range("K49").Value = 0.697
$ws.Range("L49").Value = 0.609

$ws.Range("B50").Value = "RandomForest"
$ws.Range("C50").Value = "MoCo"
$ws.Range("D50").Value = "Centers"
$ws.Range("E50").Value = "1 x 3"
$ws.Range("F50").Value = "average"
$ws.Range("M50").Value = "0.634 (0.038)"
$ws.Range("I50").Value = "standard scaling center data AFTER avg"
$ws.Range("J50").Value = 0.605
$ws.Range("K50").Value = 0.688
$ws.Range("L50").Value = 0.609

$ws.Range("B51").Value = "RandomForest"
$ws.Range("C51").Value = "MoCo"
$ws.Range("D51").Value = "Centers"
$ws.Range("E51").Value = "1 x 3"
$ws.Range("F51").Value = "average"
$ws.Range("M51").Value = "0.639 (0.025)"
$ws.Range("I51").Value = "min max scaling center data AFTER avg"
$ws.Range("J51").Value = 0.613
$ws.Range("K51").Value = 0.673
$ws.Range("L51").Value = 0.632

# ---------------------------------------------------------------------
# 2) Rename the "weakly supervision & average before predictions" label
#    (used by F5 and F13) to "weakly supervision & avg before predictions".
# ---------------------------------------------------------------------
$ws.Range("F5").Value = "weakly supervision & avg before predictions"
$ws.Range("F13").Value = "weakly supervision & avg before predictions"

# ---------------------------------------------------------------------
# 3) Grow the Tabelle1 table (and with it, its AutoFilter) down to O51.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:O51"))

# ---------------------------------------------------------------------
# 4) Extend the two row-bound conditional formats (on J and K) to O51.
# ---------------------------------------------------------------------
$cfs = $ws.Cells.FormatConditions
for ($i = 1; $i -le $cfs.Count; $i++) {
    $fc = $cfs.Item($i)
    $addr = $fc.AppliesTo.Address()
    if ($addr -eq '$K$2:$K$47') {
        $fc.ModifyAppliesToRange($ws.Range("K2:K51"))
    }
    if ($addr -eq '$J$2:$J$47') {
        $fc.ModifyAppliesToRange($ws.Range("J2:J51"))
    }
}

# ---------------------------------------------------------------------
# 5) Column width tweaks: E grows a touch, F shrinks (no longer bestFit).
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 23.998697916666668

# ---------------------------------------------------------------------
# 6) Move the view: scroll so column F / row 15 is the top-left corner,
#    and select J37 (where editing finished).
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 6
[void]$ws.Range("J37").Select()
